$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2
$ws.Range("F2").Value = 10
$ws.Range("H2").Value = 10

# Row 7
$ws.Range("E7").Value = 22
$ws.Range("F7").Value = 11
$ws.Range("H7").Value = 11

# Row 9
$ws.Range("F9").Value = 6
$ws.Range("H9").Value = 6

# Row 15
$ws.Range("E15").Value = 74
$ws.Range("F15").Value = 37
$ws.Range("H15").Value = 37

# Row 16
$ws.Range("E16").Value = 273
$ws.Range("F16").Value = 76
$ws.Range("H16").Value = 76

# Row 18
$ws.Range("E18").Value = 77
